$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# values (e.g. "0.9984") are not auto-converted to numbers by Excel.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.377.71"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.849.76"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "240.54"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").Value = "0.6336"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "0.07584"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").Value = "0.2926"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("D11").Value = "0.07745"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").Value = "1.848.54"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").Value = "5.023"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "0.6800"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "0.00001047"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").Value = "83.24"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "6.132"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").Value = "29.367.17"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").Value = "229.83"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "7.462"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "158.71"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "0.1396"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").Value = "8.455"
$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("D27").Value = "17.65"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Value = "1.418"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("D29").Value = "1.473"
$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("D30").Value = "0.05683"
$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("D31").Value = "4.123"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").Value = "1.828"
$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D34").Value = "1.155"
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("D35").Value = "0.7075"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").Value = "2.579"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").Value = "0.01828"
$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("D38").Value = "1.244.68"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("D39").Value = "2.722"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").Value = "6.420"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").Value = "0.9027"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").Value = "2.007.39"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").Value = "101.66"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "65.77"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("D46").Value = "7.137"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "0.1168"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("D48").Value = "9.034"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("D49").Value = "0.00000000116"
$ws.Range("E49").Value = "  -3.95%  "

$ws.Range("D50").Value = "0.3958"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("E51").Value = "  -0.21%  "
